$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.994324922561646
$ws.Range("B1").Value = 2.73285174369812
$ws.Range("C1").Value = 1.846433997154236
$ws.Range("D1").Value = 1.491639375686646
$ws.Range("E1").Value = 1.382882833480835
